# Fix run-estimation formula: populate EstimatedBBRun (col I) and
# EstimatedBBObserve (col J) on the "Main_Maze_Robot_Strategies" sheet,
# and drop the stale "Analysis" notes rows at the bottom. Also tidy up
# the stored selection/zoom on a couple of sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Main_Maze_Robot_Strategies")
$ws2 = $wb.Worksheets.Item("Main_Maze_Strategies")

$xlCenter = -4108

# --- Header row (row 1): center-align every used header cell A1:J1 ---
$ws1.Range("A1:J1").HorizontalAlignment = $xlCenter

# --- "Run" blocks (rows 2-21): EstimatedBBRun = BBL_Run (col E) ---
# Strategy blocks are separated by a blank row (8, 15) - copy per block
# so we never touch the blank separator rows' values.
$runBlocks = @(@(2,7), @(9,14), @(16,21))
foreach ($block in $runBlocks) {
    $first = $block[0]
    $last  = $block[1]
    for ($r = $first; $r -le $last; $r++) {
        $ws1.Cells.Item($r, 9).Value = $ws1.Cells.Item($r, 5).Value2
    }
}

# Center-align column I and J for all data/blank rows 2-21 (J stays empty).
$ws1.Range("I2:J21").HorizontalAlignment = $xlCenter

# --- "Observe" blocks (rows 23-42): EstimatedBBObserve = BBL_Observe (col D) ---
$observeBlocks = @(@(23,28), @(30,35), @(37,42))
foreach ($block in $observeBlocks) {
    $first = $block[0]
    $last  = $block[1]
    for ($r = $first; $r -le $last; $r++) {
        $ws1.Cells.Item($r, 10).Value = $ws1.Cells.Item($r, 4).Value2
    }
}

# Center-align column I and J for all data/blank rows 22-42 (I stays empty).
$ws1.Range("I22:J42").HorizontalAlignment = $xlCenter

# --- Remove the stray "Analysis" notes (rows 44-45) ---
$xlShiftUp = -4162
$ws1.Range("A44:J45").Delete($xlShiftUp) | Out-Null

# --- Update stored view state to match the authored workbook ---
$ws2.Activate()
$ws2.Range("F17").Select() | Out-Null

$ws1.Activate()
$excel.ActiveWindow.Zoom = 85
$ws1.Range("A45").Select() | Out-Null
